# Add a new row (15) to the "链表" worksheet for the palindrome
# linked-list problem, and a matching algorithm-design entry, per the
# commit "palindrome and algo design with linked list".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 15

$ws.Cells.Item($row, 1).Value = 14
$ws.Cells.Item($row, 2).Value = 234
$ws.Cells.Item($row, 3).Value = "请判断一个链表是否为回文链表。 `n例如：1-2-2-1，1-2-3-2-1"
$ws.Cells.Item($row, 4).Value = "1 找到链表的中间节点`n         链表长度是奇数，slow是链表中间节点，1-2-2-1`n         链表长度是偶数，slow是链表中间位置偏右侧的节点，1-2-3-2-1`n2 反转链表的后半部分`n         1-2 2-1`n         1-2 ，1-2-3，最后一个元素是中间节点，不用比较`n 两个链表逐个比较节点是否有相等"
$ws.Cells.Item($row, 5).Value = "链表反转`n中间节点`n链表迭代"
$ws.Cells.Item($row, 6).Value = "O(N), N是元素个数"
$ws.Cells.Item($row, 7).Value = "O(1)"

# Match the author's row height for the newly wrapped, multi-line cells.
$ws.Rows.Item($row).RowHeight = 200

# Move the on-screen selection the way the author left it after adding
# the row (one row further down than before).
$ws.Range("D22").Select() | Out-Null
